$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FE")
$tbl = $ws.ListObjects.Item(1)

# New row ("FE02") in Table1 — add first so its strings land before the
# new header / observation strings in the shared-strings table.
$newRow = $tbl.ListRows.Add()
$newRow.Range.Item(1).Value = "FE02"
$newRow.Range.Item(2).Value = 200
$newRow.Range.Item(3).Value = 9
$newRow.Range.Item(4).Value = 600
$newRow.Range.Item(5).Value = 50
$newRow.Range.Item(6).Value = 10881
$newRow.Range.Item(7).Value = 0.999
$newRow.Range.Item(8).Value = $true
$newRow.Range.Item(9).Value = $true
$newRow.Range.Item(10).Value = $true
$newRow.Range.Item(11).Value = "Todo TRUE"
$newRow.Range.Item(12).Value = "Todo TRUE"

# New column ("Observación") in Table1, with a note only on the first
# data row.
$newCol = $tbl.ListColumns.Add()
$newCol.Range.Item(1).Value = "Observación"
$ws.Range("M2").Value = "Se cortó antes de tiempo"

# Restore the on-screen selections. FE must stay the active tab, so
# touch TS first and re-select FE last.
$ts = $wb.Worksheets.Item("TS")
$ts.Range("N13").Select()

$ws.Range("E19").Select()
